# Fixed the constant Carousell bugs (hopefully) and added colour to Excel sheet.
#
# - Replaces the stale "Lego Comparison" listing rows (2-15) with fresh
#   Carousell/BrickLink data, and appends 5 new rows (16-20).
# - Highlights the "Carousell Price" column (B) using the same Red/Green/
#   Yellow ("Bad"/"Good"/"Neutral") palette Excel's built-in cell styles use,
#   to flag how each listing's asking price compares to BrickLink.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# BGR integer values (as consumed by Range.Interior.Color) for the three
# highlight colours used below.
$colorBad     = 13551615   # FFC7CE - pinkish red
$colorGood    = 13561798   # C6EFCE - green
$colorNeutral = 10284031   # FFEB9C - yellow

$data = @(
    @{ Row=2;  Title="lego creator 30688"; Price=10; Style="Bad"; Code="30688"; BLPrice=5.63 },
    @{ Row=3;  Title="lego marvel 76261"; Price=150; Style="Bad"; Code="76261"; BLPrice=99.86 },
    @{ Row=4;  Title="lego system 1252 shell tanker truck"; Price=70; Style="Bad"; Code="1252"; BLPrice=58.21 },
    @{ Row=5;  Title="lego 11033 classic 1800pcs"; Price=100; Style="Bad"; Code="11033"; BLPrice=87.06 },
    @{ Row=6;  Title="lego architecture london great britain 21034"; Price=60; Style="Bad"; Code="21034"; BLPrice=32.49 },
    @{ Row=7;  Title="lego 76294 (x-men mansion only)"; Price=140; Style="Good"; Code="76294"; BLPrice=242.84 },
    @{ Row=8;  Title="lego technic 42161 lamborghini"; Price=70; Style="Bad"; Code="42161"; BLPrice=53.58 },
    @{ Row=9;  Title="lego 60005 fire boat"; Price=35; Style="Neutral"; Code="60005"; BLPrice=72.12 },
    @{ Row=10; Title="clearance sale lego 40529 children's amusement park"; Price=14.9; Style="Bad"; Code="40529"; BLPrice=11.57 },
    @{ Row=11; Title="lego friends (41732)"; Price=220; Style="Bad"; Code="41732"; BLPrice=177.79 },
    @{ Row=12; Title="lego star wars 4486 & 4487 vintage mini building set from 2003"; Price=80; Style="Bad"; Code="4486"; BLPrice=27.18 },
    @{ Row=13; Title="lego 10729"; Price=35; Style="Neutral"; Code="10729"; BLPrice=42.47 },
    @{ Row=14; Title="lego disney castle 43205"; Price=60; Style="Good"; Code="43205"; BLPrice=110.76 },
    @{ Row=15; Title="lego 75372: clone trooper & battle droid battle pack"; Price=40; Style="Bad"; Code="75372"; BLPrice=20.26 },
    @{ Row=16; Title="lego 75337: at-te walker"; Price=180; Style="Bad"; Code="75337"; BLPrice=120.9 },
    @{ Row=17; Title="lego disney king magnifico’s castle 43224 building toy set; detailed castle makes a fun gift for ages 7 and over (613 pieces) christmas gift"; Price=60; Style="Bad"; Code="43224"; BLPrice=46.04 },
    @{ Row=18; Title="lego 76023 the tumbler- new in box"; Price=250; Style="Neutral"; Code="76023"; BLPrice=284.56 },
    @{ Row=19; Title="lego 6210 jabba’s sail barge"; Price=665; Style="Good"; Code="6210"; BLPrice=868.96 },
    @{ Row=20; Title="lego 40516 - everyone is awesome"; Price=45; Style="Bad"; Code="40516"; BLPrice=43.55 }
)

foreach ($item in $data) {
    $r = $item.Row

    $ws.Range("A$r").Value = $item.Title
    $ws.Range("B$r").Value = $item.Price

    # Set Code column holds identifiers like "30688" - keep them as text
    # (matching the rest of the sheet) instead of letting Excel coerce the
    # digit-only string into a number, then drop the helper number-format
    # style again so the cell is left on the default style.
    $cCell = $ws.Range("C$r")
    $cCell.NumberFormat = "@"
    $cCell.Value = $item.Code
    $cCell.Style = "Normal"

    $ws.Range("D$r").Value = $item.BLPrice

    switch ($item.Style) {
        "Bad"     { $ws.Range("B$r").Interior.Color = $colorBad }
        "Good"    { $ws.Range("B$r").Interior.Color = $colorGood }
        "Neutral" { $ws.Range("B$r").Interior.Color = $colorNeutral }
    }
}
